# fitur baru: statistik produksi
# Update the ATAP bulanan kab. template header labels to include their
# units of measure, matching the new "statistik produksi" feature.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Luas Panen Padi (Ha)"
$ws.Range("E1").Value = "Produksi Padi (Ton GKG)"
$ws.Range("F1").Value = "Produksi Beras (Ton)"
